$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 27, pushing existing rows 27-44 down to 28-45.
$ws.Rows.Item(27).Insert()

# Populate the newly inserted row 27 with the new weekly record.
$ws.Cells.Item(27, 1).Value  = 4
$ws.Cells.Item(27, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(27, 3).Value  = "Los Lagos"
$ws.Cells.Item(27, 4).Value  = 44488
$ws.Cells.Item(27, 5).Value  = 10
$ws.Cells.Item(27, 6).Value  = 100112026
$ws.Cells.Item(27, 7).Value  = "Haba"
$ws.Cells.Item(27, 8).Value  = "Sin especificar"
$ws.Cells.Item(27, 9).Value  = "Primera"
$ws.Cells.Item(27, 10).Value = 150
$ws.Cells.Item(27, 11).Value = 10000
$ws.Cells.Item(27, 12).Value = 10000
$ws.Cells.Item(27, 13).Value = 10000
$ws.Cells.Item(27, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(27, 15).Value = "Región Metropolitana"
$ws.Cells.Item(27, 16).Value = 400
$ws.Cells.Item(27, 17).Value = 25
$ws.Cells.Item(27, 18).Value = "Hortaliza"
